$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.912.09'
$ws.Cells.Item(2, 5).Value = '  +2.24%  '
$ws.Cells.Item(3, 4).Value = '3.051.58'
$ws.Cells.Item(3, 5).Value = '  +2.16%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '524.60'
$ws.Cells.Item(5, 5).Value = '  +5.44%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.32'
$ws.Cells.Item(6, 5).Value = '  +5.87%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  -0.02%  '
$ws.Cells.Item(8, 5).Value = '  +4.71%  '
$ws.Cells.Item(9, 5).Value = '  +4.64%  '
$ws.Cells.Item(10, 5).Value = '  +8.12%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.369'
$ws.Cells.Item(11, 5).Value = '  +5.13%  '
$ws.Cells.Item(12, 5).Value = '  +2.61%  '
$ws.Cells.Item(13, 4).Value = '3.574.65'
$ws.Cells.Item(13, 5).Value = '  +2.32%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.89'
$ws.Cells.Item(14, 5).Value = '  +8.50%  '
$ws.Cells.Item(15, 5).Value = '  +17.10%  '
$ws.Cells.Item(16, 4).Value = '57.868.43'
$ws.Cells.Item(16, 5).Value = '  +2.22%  '
$ws.Cells.Item(17, 5).Value = '  +6.90%  '
$ws.Cells.Item(18, 4).Value = '3.047.37'
$ws.Cells.Item(18, 5).Value = '  +2.20%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.92'
$ws.Cells.Item(19, 5).Value = '  +4.98%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.19'
$ws.Cells.Item(20, 5).Value = '  +5.53%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '340.30'
$ws.Cells.Item(21, 5).Value = '  +4.38%  '
$ws.Cells.Item(22, 5).Value = '  +0.02%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.499'
$ws.Cells.Item(23, 5).Value = '  +7.41%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '65.27'
$ws.Cells.Item(24, 5).Value = '  +6.50%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.172'
$ws.Cells.Item(25, 5).Value = '  +6.65%  '
$ws.Cells.Item(26, 4).Value = '0.0₃0971'
$ws.Cells.Item(26, 5).Value = '  +6.52%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.996'
$ws.Cells.Item(27, 5).Value = '  -0.04%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.00'
$ws.Cells.Item(28, 5).Value = '  +7.95%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.30'
$ws.Cells.Item(29, 5).Value = '  +7.08%  '
$ws.Cells.Item(30, 5).Value = '  +7.25%  '
$ws.Cells.Item(31, 5).Value = '  +6.40%  '
$ws.Cells.Item(32, 5).Value = '  +5.59%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '156.22'
$ws.Cells.Item(33, 5).Value = '  +1.25%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.72'
$ws.Cells.Item(34, 5).Value = '  +5.52%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.89'
$ws.Cells.Item(36, 5).Value = '  +4.00%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '25.88'
$ws.Cells.Item(37, 5).Value = '  +10.23%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0696'
$ws.Cells.Item(38, 5).Value = '  +3.15%  '
$ws.Cells.Item(39, 4).Value = '3.087.93'
$ws.Cells.Item(39, 5).Value = '  +2.30%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '37.70'
$ws.Cells.Item(40, 5).Value = '  +2.76%  '
$ws.Cells.Item(41, 5).Value = '  +8.27%  '
$ws.Cells.Item(42, 5).Value = '  +0.00%  '
$ws.Cells.Item(43, 5).Value = '  +4.71%  '
$ws.Cells.Item(44, 5).Value = '  +4.01%  '
$ws.Cells.Item(45, 4).Value = '2.329.36'
$ws.Cells.Item(45, 5).Value = '  +5.87%  '
$ws.Cells.Item(46, 5).Value = '  +2.43%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.03'
$ws.Cells.Item(47, 5).Value = '  +3.80%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0248'
$ws.Cells.Item(48, 5).Value = '  +4.42%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '6.03'
$ws.Cells.Item(49, 5).Value = '  +5.65%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '20.06'
$ws.Cells.Item(50, 5).Value = '  +4.63%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0893'
$ws.Cells.Item(51, 5).Value = '  +6.01%  '
